# Apply the cryptos.xlsx data refresh described in the commit.
# Values in column D are price strings that sometimes look numeric
# (e.g. "240.67", "0.00000000119"); the source sheet stores them as
# plain text, so we force a text number format before writing those
# cells and then restore the default "Normal" style so no stray
# formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @('D4', 'D5', 'D6', 'D8', 'D9', 'D10', 'D11', 'D12', 'D14', 'D15', 'D16', 'D17', 'D18', 'D21', 'D24', 'D26', 'D27', 'D28', 'D30', 'D31', 'D32', 'D34', 'D35', 'D37', 'D40', 'D42', 'D43', 'D44', 'D45', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '29.054.71'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').Value = '1.829.41'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('D4').Value = '0.9987'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '240.67'
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('D6').Value = '0.6229'
$ws.Range('E6').Value = '  -6.07%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '44.33'
$ws.Range('E8').Value = '  +5.63%  '
$ws.Range('D9').Value = '0.07370'
$ws.Range('E9').Value = '  -0.73%  '
$ws.Range('D10').Value = '0.2922'
$ws.Range('E10').Value = '  -0.50%  '
$ws.Range('D11').Value = '22.69'
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('D12').Value = '0.07585'
$ws.Range('E12').Value = '  -1.99%  '
$ws.Range('D13').Value = '1.829.66'
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('D14').Value = '4.962'
$ws.Range('E14').Value = '  -0.55%  '
$ws.Range('D15').Value = '0.6624'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').Value = '82.14'
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('D17').Value = '0.000009114'
$ws.Range('E17').Value = '  +8.98%  '
$ws.Range('D18').Value = '6.026'
$ws.Range('E18').Value = '  -1.19%  '
$ws.Range('D19').Value = '29.047.93'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').Value = '2.078.78'
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('D21').Value = '225.57'
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('E22').Value = '  -0.99%  '
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').Value = '7.182'
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').Value = '159.40'
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('D27').Value = '8.424'
$ws.Range('E27').Value = '  -2.31%  '
$ws.Range('D28').Value = '0.1359'
$ws.Range('E28').Value = '  -3.04%  '
$ws.Range('E29').Value = '  -0.95%  '
$ws.Range('D30').Value = '1.498'
$ws.Range('E30').Value = '  -0.96%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '4.037'
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '4.055'
$ws.Range('E32').Value = '  -1.45%  '
$ws.Range('E33').Value = '  +0.48%  '
$ws.Range('D34').Value = '0.05240'
$ws.Range('E34').Value = '  -1.50%  '
$ws.Range('D35').Value = '1.841'
$ws.Range('E35').Value = '  -1.53%  '
$ws.Range('E36').Value = '  +1.22%  '
$ws.Range('D37').Value = '0.7336'
$ws.Range('E37').Value = '  -2.61%  '
$ws.Range('E38').Value = '  +0.69%  '
$ws.Range('D39').Value = '1.285.15'
$ws.Range('E39').Value = '  +0.17%  '
$ws.Range('D40').Value = '2.747'
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('E41').Value = '  -0.82%  '
$ws.Range('D42').Value = '6.320'
$ws.Range('E42').Value = '  +6.38%  '
$ws.Range('D43').Value = '0.9022'
$ws.Range('E43').Value = '  -2.90%  '
$ws.Range('D44').Value = '1.002'
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('D45').Value = '101.75'
$ws.Range('E45').Value = '  -0.39%  '
$ws.Range('D46').Value = '1.976.86'
$ws.Range('E46').Value = '  -0.68%  '
$ws.Range('E47').Value = '  -0.67%  '
$ws.Range('D48').Value = '63.91'
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.00000000119'
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = '1.708'
$ws.Range('E50').Value = '  -3.33%  '
$ws.Range('B51').Value = 'TheSandbox'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D51').Value = '0.3966'
$ws.Range('E51').Value = '  -1.50%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

